$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.342.20"
$ws.Range("E2").Value = "  -2.88%  "
$ws.Range("D3").Value = "2.889.25"
$ws.Range("E3").Value = "  -3.94%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "583.90"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").Value = "146.08"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E8").Value = "  -2.56%  "
$ws.Range("D9").Value = "2.885.76"
$ws.Range("E9").Value = "  -4.06%  "
$ws.Range("D10").Value = "6.61"
$ws.Range("E10").Value = "  +5.96%  "
$ws.Range("E11").Value = "  -3.56%  "
$ws.Range("D12").Value = "0.445"
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").Value = "34.03"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "3.371.89"
$ws.Range("E16").Value = "  -3.77%  "
$ws.Range("D17").Value = "6.78"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "60.380.19"
$ws.Range("E18").Value = "  -2.77%  "
$ws.Range("D19").Value = "2.892.90"
$ws.Range("E19").Value = "  -3.73%  "
$ws.Range("D20").Value = "423.79"
$ws.Range("E20").Value = "  -5.12%  "
$ws.Range("D21").Value = "13.57"
$ws.Range("E21").Value = "  -4.06%  "
$ws.Range("D22").Value = "0.667"
$ws.Range("E22").Value = "  -2.68%  "
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("D24").Value = "80.69"
$ws.Range("E24").Value = "  -1.77%  "
$ws.Range("D25").Value = "11.03"
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("E26").Value = "  -3.12%  "
$ws.Range("D27").Value = "11.69"
$ws.Range("E27").Value = "  -3.40%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "7.19"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D31").Value = "2.17"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D33").Value = "26.37"
$ws.Range("E33").Value = "  -3.68%  "
$ws.Range("E34").Value = "  -4.09%  "
$ws.Range("D35").Value = "0.0₃0832"
$ws.Range("E35").Value = "  -1.77%  "
$ws.Range("E36").Value = "  -1.81%  "
$ws.Range("D37").Value = "5.65"
$ws.Range("E37").Value = "  -2.89%  "
$ws.Range("D38").Value = "49.54"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("E39").Value = "  -1.40%  "
$ws.Range("D40").Value = "2.91"
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("D41").Value = "0.122"
$ws.Range("E41").Value = "  -1.21%  "
$ws.Range("D44").Value = "41.02"
$ws.Range("E44").Value = "  -1.13%  "
$ws.Range("D47").Value = "133.02"
$ws.Range("E47").Value = "  -0.51%  "
$ws.Range("D48").Value = "2.645.07"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D50").Value = "25.11"
$ws.Range("E50").Value = "  +5.74%  "
$ws.Range("E51").Value = "  -0.91%  "

# Row 45/46 swap: VeChain moves to 45, Bittensor moves to 46
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0344"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "370.11"
$ws.Range("E46").Value = "  -5.92%  "
